$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# Values for columns C (Minute1), D (Second1), E (Rep1) for rows 2-19
$data = @{
    2  = @(15, 0, 298)
    3  = @(15, 0, 349)
    4  = @(15, 0, 300)
    5  = @(15, 0, 288)
    6  = @(15, 0, 340)
    7  = @(15, 0, 267)
    8  = @(15, 0, 317)
    9  = @(15, 0, 347)
    10 = @(15, 0, 346)
    11 = @(15, 0, 318)
    12 = @(15, 0, 313)
    13 = @(15, 0, 308)
    14 = @(15, 0, 326)
    15 = @(15, 0, 274)
    16 = @(15, 0, 296)
    17 = @(15, 0, 324)
    18 = @(15, 0, 332)
    19 = @(15, 0, 287)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
}

# Update the active cell selection to E20
$ws.Range("E20").Select()
